# Apply the changes described by the diff:
# Sheet "Metadata" (first sheet):
#   B7  ("Experimental" row) : "" -> "false"
#   B8  ("Date" row)         : "2025-11-28T14:35:57+00:00" -> "2025-11-30T13:08:37+00:00"
#   B17 ("Description" row)  : "" -> "Methods for determining maximum heart rate"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row value: set as literal text "false" (not a native boolean),
# matching how this FHIR IG export tool stores boolean-like metadata as text
# (see row 20 "Immutable" = "BooleanType[null]"). A leading apostrophe forces
# Excel to store it as text instead of auto-converting to a TRUE/FALSE boolean.
$ws.Range("B7").Value = "'false"

# "Date" row value update
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# "Description" row value (was previously blank)
$ws.Range("B17").Value = "Methods for determining maximum heart rate"
